# Updates cryptos list: Price (D) and Volume(1h) (E) columns for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look numeric need to be forced to stay as
# plain text (matching the workbook's existing inline-string "Price" column),
# otherwise Excel auto-converts them to numbers. We set NumberFormat = "@" first,
# then revert the cell style back to Normal afterwards so no stray formatting
# is introduced.
$textRows = 4,5,6,7,8,11,13,14,17,19,20,21,22,24,25,27,28,29,32,33,35,36,37,38,39,40,41,42,43,44,45,46,47,50,51
foreach ($r in $textRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = "56.967.18"
$ws.Range("E2").Value = "  +1.59%  "
$ws.Range("D3").Value = "2.345.08"
$ws.Range("E3").Value = "  +1.37%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.61%  "
$ws.Range("D5").Value = "518.27"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").Value = "136.06"
$ws.Range("E6").Value = "  +2.16%  "
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "0.538"
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").Value = "2.357.70"
$ws.Range("E9").Value = "  +0.66%  "
$ws.Range("E10").Value = "  -0.80%  "
$ws.Range("D11").Value = "5.42"
$ws.Range("E11").Value = "  +5.15%  "
$ws.Range("E12").Value = "  -1.61%  "
$ws.Range("D13").Value = "0.344"
$ws.Range("E13").Value = "  -0.14%  "
$ws.Range("D14").Value = "24.02"
$ws.Range("E14").Value = "  -0.83%  "
$ws.Range("D15").Value = "2.754.95"
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("D16").Value = "56.913.63"
$ws.Range("E16").Value = "  +1.26%  "
$ws.Range("D17").Value = "0.0000135"
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("D18").Value = "2.344.89"
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("D19").Value = "10.60"
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").Value = "327.39"
$ws.Range("E20").Value = "  +1.49%  "
$ws.Range("D21").Value = "4.23"
$ws.Range("E21").Value = "  -0.44%  "
$ws.Range("D22").Value = "6.78"
$ws.Range("E22").Value = "  +2.08%  "
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").Value = "61.00"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").Value = "0.166"
$ws.Range("E25").Value = "  +4.56%  "
$ws.Range("E26").Value = "  +0.95%  "
$ws.Range("D27").Value = "7.96"
$ws.Range("E27").Value = "  +3.64%  "
$ws.Range("D28").Value = "1.31"
$ws.Range("E28").Value = "  +10.25%  "
$ws.Range("D29").Value = "170.33"
$ws.Range("E29").Value = "  -1.26%  "
$ws.Range("D30").Value = "0.0₃0745"
$ws.Range("E30").Value = "  +2.37%  "
$ws.Range("E31").Value = "  +1.41%  "
$ws.Range("D32").Value = "6.25"
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("D33").Value = "18.60"
$ws.Range("E33").Value = "  +1.17%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").Value = "0.995"
$ws.Range("E35").Value = "  +0.27%  "
$ws.Range("D36").Value = "1.27"
$ws.Range("E36").Value = "  +0.61%  "
$ws.Range("D37").Value = "0.918"
$ws.Range("E37").Value = "  -1.35%  "
$ws.Range("D38").Value = "4.03"
$ws.Range("E38").Value = "  +1.06%  "
$ws.Range("D39").Value = "1.57"
$ws.Range("E39").Value = "  +3.40%  "
$ws.Range("D40").Value = "38.45"
$ws.Range("E40").Value = "  +2.75%  "
$ws.Range("D41").Value = "148.34"
$ws.Range("E41").Value = "  +7.15%  "
$ws.Range("D42").Value = "0.383"
$ws.Range("E42").Value = "  -0.57%  "
$ws.Range("D43").Value = "3.63"
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("D44").Value = "282.53"
$ws.Range("E44").Value = "  +5.92%  "
$ws.Range("D45").Value = "5.27"
$ws.Range("E45").Value = "  +2.47%  "
$ws.Range("D46").Value = "0.0938"
$ws.Range("E46").Value = "  +1.43%  "
$ws.Range("D47").Value = "0.0506"
$ws.Range("E47").Value = "  -1.06%  "
$ws.Range("E48").Value = "  +0.87%  "
$ws.Range("E49").Value = "  +1.70%  "
$ws.Range("D50").Value = "18.11"
$ws.Range("E50").Value = "  +6.75%  "
$ws.Range("D51").Value = "17.58"
$ws.Range("E51").Value = "  +3.44%  "

# Revert the temporary NumberFormat tweak so cell styling matches the original
# (Normal/General) while the values remain text.
foreach ($r in $textRows) {
    $ws.Range("D$r").Style = "Normal"
}

